$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently has 190 data rows (rows 2-191), then an empty separator
# row (192) and a filter-description row (193). The updated export contains
# 211 data rows, so insert 21 additional blank rows before the separator row
# to make room, pushing the separator/filter rows down to 213/214.
$ws.Range("A192:A212").EntireRow.Insert() | Out-Null

# Account numbers ("Conta") are 9-digit codes with leading zeros and must be
# stored as text, not numbers, otherwise Excel would strip the leading zeros.
$ws.Range("A2:A212").NumberFormat = "@"

$rows = @(
  @("005554830", "PAULO", 66178.74),
  @("004222784", "RAFAEL", 55148.95),
  @("005079311", "JOVINO", 50737.03),
  @("004419765", "WALTER", 46325.12),
  @("005639338", "CLOTILDE", 22350.19),
  @("004368468", "AHMAD", 21621.63),
  @("004381194", "ALINNE", 16577.85),
  @("004473829", "LUCAS", 16544.69),
  @("004454491", "MARCO", 14338.73),
  @("004994036", "BALTASAR", 11560.28),
  @("004693349", "CATARINE", 11010.38),
  @("004497825", "PRISCILLA", 10003.11),
  @("004267044", "PATRICIA", 7720.86),
  @("005040864", "ANDRE", 7000),
  @("004584982", "BRUNO", 6617.88),
  @("004243043", "SUELI", 6612.04),
  @("002064834", "RAFAELA", 5514.9),
  @("003497496", "ELISANDRA", 5514.9),
  @("004508516", "EDUARDO", 5514.9),
  @("004500087", "JANNE", 5505.18),
  @("004582648", "ANUAR", 5000),
  @("005547703", "SILVIA", 4427.61),
  @("004313254", "GUSTAVO", 4292),
  @("005890232", "TAYLA", 4000),
  @("004472538", "RODOLFO", 3308.93),
  @("004500804", "RAFAEL", 3308.93),
  @("004643153", "CARLA", 3308.93),
  @("004713953", "ALESSANDRA", 1170.04),
  @("004355790", "MINEIA", 993.05),
  @("004392159", "RODRIGO", 900.21),
  @("003435941", "HEITOR", 807.5),
  @("005646524", "EVANGELINA", 785.39),
  @("004499920", "FABIANO", 730.93),
  @("004839302", "GIANLUCA", 530.05),
  @("004479287", "ANA", 497.82),
  @("005002457", "ROSANGELA", 379.59),
  @("004376853", "ALBERTO", 260.19),
  @("004165558", "GABRIEL", 200),
  @("004436055", "MARCO", 124.15),
  @("004467884", "ANA", 100),
  @("001368670", "THIAGO", 99.84),
  @("004211911", "ZENILDA", 96),
  @("004431591", "MARIO", 94.24),
  @("005314853", "GLEUBER", 93.89),
  @("005681354", "MATHEUS", 93.63),
  @("005266369", "EG", 93.11),
  @("005696533", "CARLOS", 90.18),
  @("004546050", "LUIS", 88.61),
  @("004212132", "JOAO", 86.38),
  @("004360431", "CARLOS", 86.31),
  @("004809902", "PEDRO", 85.9),
  @("004877741", "LUIZ", 85.27),
  @("004206790", "EMMANUELLE", 84.05),
  @("004268684", "PATRICIA", 84.05),
  @("004482163", "CIRIA", 84.05),
  @("005079458", "JONAS", 83.23),
  @("004480134", "JOSE", 82.33),
  @("008030888", "SONIA", 81.06),
  @("004424761", "PEDRO", 80),
  @("004754920", "LUIS", 79.23),
  @("005685089", "CARNEIRO", 78.85),
  @("005981575", "GLAUCIANE", 77.19),
  @("005142611", "GUILHERME", 74.22),
  @("001719494", "LUIS", 73.48),
  @("005558076", "ALEXANDRE", 72.87),
  @("004211922", "CARLOS", 72.84),
  @("004342617", "JURACI", 72.84),
  @("005428871", "ROSANGELA", 72.84),
  @("005761337", "MICHELE", 72.59),
  @("004563252", "FERNANDO", 70.58),
  @("005902068", "VIVIANE", 69.94),
  @("005077648", "DUNAS", 69.8),
  @("004908680", "ELENE", 69.7),
  @("004643880", "GABRIEL", 69.55),
  @("004462930", "WALTER", 69.32),
  @("004924222", "ROSANE", 68.85),
  @("002823185", "SIMONE", 68.42),
  @("005007114", "STELLA", 67.41),
  @("004216434", "JAIME", 67.16),
  @("004457389", "RAFAEL", 66.93),
  @("004332207", "IRACY", 66.88),
  @("004749928", "NADY", 66.86),
  @("008007759", "CRISTINA", 66.53),
  @("004204242", "MARIA", 66.15),
  @("004550605", "REJANE", 65.71),
  @("004756968", "DANIELY", 63.91),
  @("004277637", "LARA", 63.37),
  @("004335251", "EDMUNDO", 62.82),
  @("004472760", "SANDRA", 62.46),
  @("004374891", "RODRIGO", 61.18),
  @("004384131", "ANDRE", 58.71),
  @("004382374", "THEOMAR", 57.75),
  @("004974089", "CELIA", 57.68),
  @("004207955", "SILVANIA", 56.03),
  @("004212409", "RAFAEL", 56.03),
  @("004216657", "JOAO", 56.03),
  @("004455489", "SOLANGE", 56.03),
  @("005046919", "MARIANA", 56.03),
  @("005061124", "BRUNO", 56.03),
  @("005105970", "VERA", 56.03),
  @("004552021", "MARIA", 53.8),
  @("004224815", "GUILHERME", 52.26),
  @("004400640", "FELIPE", 51.44),
  @("005076418", "LEONARDO", 48.27),
  @("004508159", "FELIPE", 46.28),
  @("004517080", "TATIANA", 45.62),
  @("004340984", "RENATA", 45.24),
  @("001731007", "GUILHERME", 44.73),
  @("005684392", "PAULO", 44.25),
  @("000827730", "LUCIANA", 41.09),
  @("004886366", "RENATO", 40.73),
  @("005295509", "BHRUNA", 40.55),
  @("004290978", "LARISSA", 40.34),
  @("004971448", "CLOVIS", 40.07),
  @("005245032", "ROSA", 39.91),
  @("004238164", "DANIELA", 39.09),
  @("005715733", "ADRIANA", 38.93),
  @("004260002", "ERICA", 38.6),
  @("004752519", "MARCUS", 37.58),
  @("005173958", "VENIA", 37.39),
  @("004413537", "CLAUDIA", 36.55),
  @("004806286", "VERA", 35.77),
  @("004240292", "MARCO", 34.71),
  @("004218542", "JOSE", 31.67),
  @("004207374", "ANGELICA", 31.25),
  @("004481463", "MARA", 30.77),
  @("005927101", "SIMONE", 30),
  @("004027477", "GABRIELA", 28.97),
  @("004377415", "ANGELA", 28.73),
  @("001294033", "VIVIANE", 28.72),
  @("004230529", "LAIS", 28.18),
  @("005305965", "SIDMAR", 28.12),
  @("004404724", "LEANDRO", 26.71),
  @("004350197", "GISELA", 25.08),
  @("004691225", "ANNA", 24.56),
  @("004813088", "JULIANA", 23.58),
  @("004459461", "INTERLAGOS", 23.46),
  @("008070544", "MARINA", 22.74),
  @("008090243", "GABRIEL", 22.74),
  @("004371857", "NAZARETH", 21.52),
  @("004388077", "WLADMIR", 20.89),
  @("004214604", "MARIA", 20.75),
  @("005143579", "GABRIEL", 19.2),
  @("004204255", "AMADO", 18.77),
  @("002894447", "JOAO", 16.48),
  @("004332103", "JOSE", 16.02),
  @("004422594", "WANDIR", 14.67),
  @("005878792", "JUNIO", 14.49),
  @("004920447", "MARILIA", 13.59),
  @("004454365", "RAFAEL", 13.38),
  @("004264780", "MARCELO", 12.67),
  @("001804114", "WAGNER", 12.43),
  @("004374943", "LEONARDO", 12.08),
  @("004472076", "RUBENS", 11.8),
  @("004335031", "EDMUNDO", 11.17),
  @("005374916", "MARCO", 10.89),
  @("004216298", "FLORDELIZ", 9.79),
  @("008013889", "CAROLINA", 9.72),
  @("005142592", "ALBERTO", 8.87),
  @("004308815", "ZELI", 7.54),
  @("004505474", "RICARDO", 7.46),
  @("004530494", "ROSANGELA", 6.99),
  @("004425965", "CAROLLINA", 6.9),
  @("004907688", "HEITOR", 6.12),
  @("004224405", "MILA", 5.82),
  @("004278212", "LEONARDO", 5.55),
  @("004448501", "JOAO", 5.55),
  @("005268516", "LUIS", 4.91),
  @("004381415", "JOAO", 4.46),
  @("000834301", "MARCUS", 4.4),
  @("008012870", "ANA", 4.3),
  @("008004995", "JOSE", 3.74),
  @("004220849", "DULCE", 3.59),
  @("004431546", "GABRIELA", 2.9),
  @("004451652", "MATEUS", 2.87),
  @("005022526", "ALEXANDRE", 1.7),
  @("000431814", "GUILHERME", 1.1),
  @("004360430", "VIOMAR", 1),
  @("004486497", "ELENA", 0.96),
  @("004115403", "HEBERT", 0.88),
  @("005660155", "CAROLINA", 0.85),
  @("004223502", "BRUNA", 0.78),
  @("004587511", "CARLOS", 0.69),
  @("004380948", "LUISA", 0.62),
  @("004473942", "DAIANNE", 0.62),
  @("003894173", "ANDREA", 0.48),
  @("004453302", "ISABELLA", 0.39),
  @("004638738", "GABRIEL", 0.29),
  @("004278033", "DAISY", 0.21),
  @("005662526", "AGUINALDO", 0.18),
  @("002687737", "JOSE", 0.17),
  @("004432455", "LUCIANA", 0.17),
  @("002694089", "VITOR", 0.16),
  @("004357159", "JOAO", 0.15),
  @("004320840", "NATALIA", 0.14),
  @("000772433", "MARCELO", 0.13),
  @("001000288", "ISABELLA", 0.13),
  @("005530256", "CAROLINA", 0.1),
  @("004451996", "ADRIANO", 0.09),
  @("005047946", "GABRIEL", 0.09),
  @("004223226", "YESHUA", 0.04),
  @("004281300", "FRANKLIN", 0.04),
  @("005274028", "RAFAEL", 0.04),
  @("004329229", "GABRIEL", 0.03),
  @("004213373", "ALEXANDRE", 0.02),
  @("004339183", "JALISON", 0.02),
  @("004870976", "HFR", 0.02),
  @("000938440", "BASE", 0.01),
  @("002878817", "GUILHERME", 0.01),
  @("004400000", "VILMA", 0.01),
  @("004976625", "NORTON", 0.01)
)

$n = $rows.Count
$data = New-Object 'object[,]' $n,3
for ($i = 0; $i -lt $n; $i++) {
    $data[$i,0] = $rows[$i][0]
    $data[$i,1] = $rows[$i][1]
    $data[$i,2] = $rows[$i][2]
}

$ws.Range("A2:C212").Value = $data
